$d = $word.ActiveDocument

# Collapse to the very end of the document body (after the last existing
# paragraph, before the sectPr) so the new content is appended cleanly.
$ip = $d.Content
$ip.Collapse(0)
$insertStart = $ip.Start

# Build the new paragraphs as a raw OOXML "flat" package fragment. This lets
# us place <w:proofErr/> markers exactly where Word would (splitting runs
# around "welcome" / "github") and keep "last-minute" as its own run, which
# plain Range.InsertAfter() text insertion cannot reproduce (adjoining runs
# with identical formatting get silently merged back together).
#
# The two hyperlink URLs are inserted here as plain text; the Hyperlink
# character style + actual hyperlink relationship is applied afterwards via
# Hyperlinks.Add (see below) because InsertXML does not round-trip rStyle
# references placed directly in the fragment's rPr.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">It is totally optional, but you are also </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>welcome</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve"> to peek at the PowerPoint slides for the lecture in advance. </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">These are close to final, but I reserve the right to make small </w:t>
            </w:r>
            <w:r>
              <w:t>last-minute</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> changes. </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve">The slides are located on my </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>github</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> site as a PowerPoint file or a PDF file with speaker notes.</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>https://github.com/pmean/papers-and-presentations/blob/master/dark-side/2022-talk.pptx</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>https://github.com/pmean/papers-and-presentations/blob/master/dark-side/2022-talk-speaker-notes.pdf</w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
          </w:p>
          <w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$ip.InsertXML($xml) | Out-Null

# Turn the two plain-text URLs into real hyperlinks. Hyperlinks.Add() both
# wraps the text in <w:hyperlink> with the correct relationship id and
# applies the "Hyperlink" character style (<w:rStyle w:val="Hyperlink"/>)
# exactly like the other hyperlinks already in this document - and unlike
# directly setting Range.Style, it does not stamp stray w:rsidP attributes
# onto unrelated paragraphs.
$url1 = "https://github.com/pmean/papers-and-presentations/blob/master/dark-side/2022-talk.pptx"
$url2 = "https://github.com/pmean/papers-and-presentations/blob/master/dark-side/2022-talk-speaker-notes.pdf"

$pCount = $d.Paragraphs.Count
$hyperlinkPara1 = $d.Paragraphs.Item($pCount - 2)
$hyperlinkPara2 = $d.Paragraphs.Item($pCount - 1)

$linkRange1 = $d.Range($hyperlinkPara1.Range.Start, $hyperlinkPara1.Range.Start + $url1.Length)
$d.Hyperlinks.Add($linkRange1, $url1) | Out-Null

$linkRange2 = $d.Range($hyperlinkPara2.Range.Start, $hyperlinkPara2.Range.Start + $url2.Length)
$d.Hyperlinks.Add($linkRange2, $url2) | Out-Null
